$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test_schedule_1")

$ws.Range("C2").Value = "test_files/vids/test.mp4"
$ws.Range("C3").Value = "test_files/vids/test2.mp4"
$ws.Range("C4").Value = "test_files/vids/test4.mp4"
